$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top Gainers")

$rows = @(
  @{Row=2; B="SOLARWORLD"; C=14.1379; D=10.1482; E=5.6878},
  @{Row=3; B="BLUEDART"; C=14.0767; D=13.5028; E=10.7663},
  @{Row=4; B="DREDGECORP"; C=12.0259; D=16.5263; E=17.2881},
  @{Row=5; B="ADANIGREEN"; C=11.432; D=8.524900000000001; E=8.974},
  @{Row=6; B="INFOBEAN"; C=9.9924; D=23.0892; E=38.1616},
  @{Row=7; B="VBL"; C=9.424200000000001; D=7.7048; E=12.0014},
  @{Row=8; B="BUTTERFLY"; C=9.023999999999999; D=11.972; E=14.6049},
  @{Row=9; B="HEG"; C=8.059100000000001; D=12.0613; E=14.4583},
  @{Row=10; B="RPOWER"; C=7.326; D=3.7168; E=5.6808},
  @{Row=11; B="WALCHANNAG"; C=7.1578; D=4.683; E=-4.5619},
  @{Row=12; B="UTKARSHBNK"; C=6.9815; D=-4.0074; E=-0.6673},
  @{Row=13; B="ABREL"; C=6.8217; D=7.5963; E=7.1505},
  @{Row=14; B="SAIL"; C=6.5754; D=8.7981; E=4.7368},
  @{Row=15; B="M&MFIN"; C=6.5355; D=7.1608; E=16.0552},
  @{Row=16; B="JISLJALEQS"; C=6.3847; D=5.6756; E=-0.4377},
  @{Row=17; B="CELLO"; C=6.2926; D=5.1213; E=15.1042},
  @{Row=18; B="SANDUMA"; C=6.2666; D=3.7749; E=32.366},
  @{Row=19; B="FISCHER"; C=6.0598; D=11.0589; E=4.2252},
  @{Row=20; B="ADANIENSOL"; C=5.8673; D=3.283; E=11.8085},
  @{Row=21; B="POKARNA"; C=5.7671; D=-1.6492; E=18.5397},
  @{Row=22; B="GRAPHITE"; C=5.692; D=12.1215; E=12.3233},
  @{Row=23; B="VAIBHAVGBL"; C=5.4518; D=5.9775; E=12.2997},
  @{Row=24; B="IOC"; C=5.4232; D=8.332800000000001; E=8.7523},
  @{Row=25; B="ABDL"; C=5.3767; D=4.2775; E=27.0041},
  @{Row=26; B="EPACKPEB"; C=5.2943; D=-1.2138; E="N/A"},
  @{Row=27; B="STLTECH"; C=5.1526; D=2.1317; E=8.2811},
  @{Row=28; B="ATGL"; C=5.0242; D=4.7627; E=4.3019},
  @{Row=29; B="MEGASOFT"; C=4.9974; D=15.7588; E=33.5271},
  @{Row=30; B="PROZONER"; C=4.9921; D=15.7468; E=36.095},
  @{Row=31; B="STALLION"; C=4.9914; D=-5.2229; E=21.4391},
  @{Row=32; B="INDOTHAI"; C=4.9883; D=4.7163; E=43.9974},
  @{Row=33; B="SURYAROSNI"; C=4.9297; D=11.3293; E=2.9689},
  @{Row=34; B="HITECHGEAR"; C=4.8651; D=2.1287; E=10.9905},
  @{Row=35; B="SGMART"; C=4.7391; D=8.742100000000001; E=2.9958},
  @{Row=36; B="BAJAJINDEF"; C=4.7085; D=3.6272; E=10.6547},
  @{Row=37; B="GMBREW"; C=4.6806; D=0.2175; E=80.37439999999999},
  @{Row=38; B="AXISCADES"; C=4.6582; D=7.128; E=-2.8858},
  @{Row=39; B="DATAMATICS"; C=4.5164; D=6.9223; E=15.3061},
  @{Row=40; B="SAMBHV"; C=4.5075; D=2.9912; E=5.5433},
  @{Row=41; B="CMSINFO"; C=4.4372; D=3.2086; E=3.4159},
  @{Row=42; B="GENUSPOWER"; C=4.3919; D=2.709; E=-0.3387},
  @{Row=43; B="MRPL"; C=4.3623; D=9.813499999999999; E=20.1672},
  @{Row=44; B="GPIL"; C=4.3499; D=6.5242; E=14.6528},
  @{Row=45; B="PROSTARM"; C=4.329; D=1.4352; E=-7.5334},
  @{Row=46; B="JKIL"; C=4.2497; D=3.0575; E=1.8683},
  @{Row=47; B="TMB"; C=4.118; D=7.8418; E=15.0492},
  @{Row=48; B="SUNFLAG"; C=4.1106; D=4.447; E=4.7456},
  @{Row=49; B="FIVESTAR"; C=4.0965; D=4.1061; E=4.1838},
  @{Row=50; B="NBCC"; C=3.9797; D=2.695; E=7.1162},
  @{Row=51; B="STAR"; C=3.9496; D=3.8794; E=3.1136},
  @{Row=52; B="DCMSHRIRAM"; C=3.8931; D=10.476; E=17.8674},
  @{Row=53; B="LLOYDSENT"; C=3.8867; D=1.1737; E=10.5128},
  @{Row=54; B="SRM"; C=3.8659; D=3.5572; E=4.4692},
  @{Row=55; B="TCI"; C=3.803; D=3.7068; E=4.2072},
  @{Row=56; B="RAJRATAN"; C=3.795; D=1.3359; E=27.4528},
  @{Row=57; B="MAITHANALL"; C=3.7917; D=2.8312; E=2.1015},
  @{Row=58; B="HCC"; C=3.7446; D=2.6464; E=7.3828},
  @{Row=59; B="GAIL"; C=3.7207; D=2.2539; E=4.9975},
  @{Row=60; B="SHK"; C=3.6851; D=2.4377; E=-1.8843},
  @{Row=61; B="SUNDROP"; C=3.6802; D=2.1681; E=0.2812},
  @{Row=62; B="GPPL"; C=3.6795; D=2.6785; E=4.3093},
  @{Row=63; B="VINCOFE"; C=3.6605; D=10.5195; E=8.895799999999999},
  @{Row=64; B="RECLTD"; C=3.6455; D=2.6302; E=2.5613},
  @{Row=65; B="DCW"; C=3.643; D=2.2121; E=-4.0783},
  @{Row=66; B="ORIENTTECH"; C=3.6279; D=0.332; E=32.424},
  @{Row=67; B="INDORAMA"; C=3.6268; D=2.7915; E=13.9312},
  @{Row=68; B="ICRA"; C=3.6236; D=4.3033; E=2.7095},
  @{Row=69; B="SUZLON"; C=3.5752; D=8.194000000000001; E=5.7766},
  @{Row=70; B="ASHAPURMIN"; C=3.4856; D=6.134; E=1.9278},
  @{Row=71; B="PRAKASH"; C=3.4167; D=4.3192; E=1.0705},
  @{Row=72; B="AVALON"; C=3.4124; D=8.174200000000001; E=20.1209},
  @{Row=73; B="MSPL"; C=3.4071; D=2.0649; E=-5.2055},
  @{Row=74; B="RHIM"; C=3.3974; D=2.9716; E=4.9218},
  @{Row=75; B="BLACKBUCK"; C=3.39; D=2.0099; E=7.9628},
  @{Row=76; B="SALASAR"; C=3.372; D=4.3617; E=10.5975}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}